$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 212, shifting existing rows 212:231 down to 213:232
$ws.Rows.Item(212).Insert()

# Populate the new weekly record in row 212
$ws.Range("A212").Value = 5
$ws.Range("B212").Value = "Macroferia Regional de Talca"
$ws.Range("C212").Value = "Maule"
$ws.Range("D212").Value = 44578
$ws.Range("E212").Value = 7
$ws.Range("F212").Value = 100112003
$ws.Range("G212").Value = "Ajo"
$ws.Range("H212").Value = "Chino"
$ws.Range("I212").Value = "Primera"
$ws.Range("J212").Value = 230
$ws.Range("K212").Value = 20000
$ws.Range("L212").Value = 20000
$ws.Range("M212").Value = 20000
$ws.Range("N212").Value = "$/caja 10 kilos"
$ws.Range("O212").Value = "China"
$ws.Range("P212").Value = 2000
$ws.Range("Q212").Value = 10
$ws.Range("R212").Value = "Hortaliza"
